$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-13 Thursday", 2) | Out-Null
$d.Content.Find.Execute("991÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷6=", 2) | Out-Null
$d.Content.Find.Execute("382÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "173÷2=", 2) | Out-Null
$d.Content.Find.Execute("173÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "570÷8=", 2) | Out-Null
$d.Content.Find.Execute("168÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "753÷4=", 2) | Out-Null
$d.Content.Find.Execute("886÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "330÷3=", 2) | Out-Null
$d.Content.Find.Execute("458÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "956÷4=", 2) | Out-Null
$d.Content.Find.Execute("676÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "504÷7=", 2) | Out-Null
$d.Content.Find.Execute("271÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "841÷4=", 2) | Out-Null
$d.Content.Find.Execute("788÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "300÷5=", 2) | Out-Null
$d.Content.Find.Execute("942÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "632÷6=", 2) | Out-Null
$d.Content.Find.Execute("120÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "327÷8=", 2) | Out-Null
$d.Content.Find.Execute("586÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "668÷9=", 2) | Out-Null
$d.Content.Find.Execute("704÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "207÷2=", 2) | Out-Null
$d.Content.Find.Execute("300÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "515÷8=", 2) | Out-Null
$d.Content.Find.Execute("192÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "319÷3=", 2) | Out-Null
$d.Content.Find.Execute("611÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "348÷5=", 2) | Out-Null
$d.Content.Find.Execute("101÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "882÷3=", 2) | Out-Null
$d.Content.Find.Execute("522÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "445÷7=", 2) | Out-Null
$d.Content.Find.Execute("455÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "307÷6=", 2) | Out-Null
$d.Content.Find.Execute("100÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "857÷4=", 2) | Out-Null
$d.Content.Find.Execute("175÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷5=", 2) | Out-Null
$d.Content.Find.Execute("919÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "581÷8=", 2) | Out-Null
$d.Content.Find.Execute("469÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "691÷4=", 2) | Out-Null
$d.Content.Find.Execute("343÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "261÷7=", 2) | Out-Null
$d.Content.Find.Execute("825÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "675÷8=", 2) | Out-Null
